# Apply the edits described in the commit:
#  1. Rename "Sheet2" -> "인증리스트" (the defined name _FilterDatabase that
#     references Sheet2 updates automatically because it tracks the sheet
#     by name).
#  2. Append a new certification-list row (row 1110) to the MASTER sheet
#     with the ARTESYN PSU (700-014464-0100) entry.
#  3. Update the MASTER sheet's active selection / frozen-pane scroll
#     position to reflect the newly added row, and nudge the 인증리스트
#     sheet's scroll position as well.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename Sheet2 to 인증리스트
# ---------------------------------------------------------------------
$wsCert = $wb.Worksheets.Item("Sheet2")
$wsCert.Name = "인증리스트"

# ---------------------------------------------------------------------
# 2) Add the new row (1110) to MASTER
# ---------------------------------------------------------------------
$wsMaster = $wb.Worksheets.Item("MASTER")

# Column B ("HS Code") holds a digit-only string that must stay TEXT
# (matches the existing "quote-prefixed" numeric-text cells like B1107/
# B1108 elsewhere in the sheet) rather than being interpreted as a number.
$wsMaster.Range("B1110").NumberFormat = "@"
$wsMaster.Range("B1110").Value = "8504409011"
$wsMaster.Range("B1110").NumberFormat = "General"

# The remaining text columns are written in the same order the brand-new
# shared-string table entries appear in (C, E, G, I, M, then D, F), so the
# new <si> entries land at the same indices as the authored workbook.
$wsMaster.Range("C1110").Value = "MSF-064849 (700-014464-0100)"
$wsMaster.Range("E1110").Value = "700-014464-0100"
$wsMaster.Range("G1110").Value = "YU10710-17001A"
$wsMaster.Range("I1110").Value = "(200 - 240) V~, (50 / 60) Hz, 3.5 A (Output : 12.25 Vd.c., 134.7 A)"
$wsMaster.Range("M1110").Value = "Zhongshan Artesyn Technologies.Co.,Ltd"
$wsMaster.Range("D1110").Value = "ARTESYN - PSU,ARTESYN,PS1650,N MODEL:700-014464-0100 VOLT:200-240V"
$wsMaster.Range("F1110").Value = "MSIP-REM-AIL-700-014464"

$wsMaster.Range("H1110").Value = "화학"
$wsMaster.Range("L1110").Value = "컴퓨터용 전원공급장치"

# Re-apply the formatting used by the other rows in this block: column B
# keeps the quote-prefixed "number as text" look, column I keeps the
# rated-voltage column's centered-alignment style.
$wsMaster.Range("B1107").Copy()
$wsMaster.Range("B1110").PasteSpecial(-4122)
$wsMaster.Range("I9").Copy()
$wsMaster.Range("I1110").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) Update view/selection state
# ---------------------------------------------------------------------
$wsMaster.Activate() | Out-Null
$wsMaster.Range("C1114").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1091

$wsCert.Activate() | Out-Null
$wsCert.Range("A166:A187").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 168

$wsMaster.Activate() | Out-Null
